$wb = $excel.ActiveWorkbook

# --- "survey" sheet: drop the begin/end group+table scratch rows (11-21)
# and replace them with a single "select"/"Table" row (row 10). ---
$ws = $wb.Worksheets.Item("survey")
$ws.Range("A11:F21").EntireRow.Delete()

$ws.Range("A10").Value = "text"
$ws.Range("B10").Value = "select"
$ws.Range("C10").Value = "Table"
$ws.Range("B10").Select()

# --- "settings" sheet: fix up the form_id value and make this the
# active/selected sheet+cell. ---
$ws3 = $wb.Worksheets.Item("settings")
$ws3.Range("B2").Value = "Justtest"
$ws3.Activate()
$ws3.Range("B2").Select()
